$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert both blank rows first, from the bottom up so earlier inserts
# don't shift the row numbers we still need to target.

# New row for "MAW_FAMES" goes directly above the existing "MBW_FAMES"
# row (currently row 6).
$ws.Rows.Item(6).Insert()

# New row for "MAX_FAMES" goes directly above the existing "MBX_FAMES"
# row (currently row 2).
$ws.Rows.Item(2).Insert()

# A freshly-inserted row inherits the bold header style from row 1.
# Re-stamp the font explicitly so the new row ends up with its own
# (non-header, non-bold) style rather than staying flagged as bold.
$ws.Range("A2:E2").Font.Name = "Calibri"

# Now fill in the new rows top-to-bottom so new shared-string entries
# are appended in the same order as in the target workbook.
$ws.Cells.Item(2,1).Value = "MAX_FAMES"
$ws.Cells.Item(2,2).Value = "top soil"
$ws.Cells.Item(2,3).Value = $false
$ws.Cells.Item(2,4).Value = 0
$ws.Cells.Item(2,5).Value = $false

$ws.Cells.Item(7,1).Value = "MAW_FAMES"
$ws.Cells.Item(7,2).Value = "top soil"
$ws.Cells.Item(7,3).Value = $false
$ws.Cells.Item(7,4).Value = 0
$ws.Cells.Item(7,5).Value = $true

$ws.Range("F12").Select() | Out-Null
